# Re-shuffle the per-trial data (category/condition/stimulus/ratings) across rows
# 2..41 of the active sheet, per the commit: "make only 20 different versions and
# duplicate many times for 1000 subjects". Columns A-G (subject_id/task/block.../
# target_cat) and J (cond_mem, always blank) stay put; columns H,I,K,L,M,N,O,P,Q,R,S,T,U,V
# (category, cond_cat, correct_answer, stimulus, conceptual, perceptual, typicality,
# n, p_typicality, p_conceptual, p_perceptual, r_typicality, r_conceptual, r_perceptual)
# get redistributed across rows according to a fixed permutation of source rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 41

# destination row (index into this array is destRow-2) -> source row to copy the
# H..V (minus J) payload from.
$sourceRowFor = @(26,27,24,20,35,18,4,9,8,34,36,14,29,37,40,16,38,11,28,3,30,32,12,23,22,39,13,17,19,5,10,7,21,33,6,2,15,41,31,25)

# Columns (1-based) whose values travel together as one trial's payload.
$cols = @(8,9,11,12,13,14,15,16,17,18,19,20,21,22)

# Snapshot every payload column for every row first -- the mapping is an arbitrary
# permutation (not a simple shift), so rows get overwritten out of order and we must
# not read a row after it has already been rewritten.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

for ($i = 0; $i -lt $sourceRowFor.Length; $i++) {
    $destRow = $firstRow + $i
    $srcRow = $sourceRowFor[$i]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
